{"js": "// Word Weekly Progress Report heading: bump \"#3\" to \"#4\".\n// Word keeps a \"_GoBack\" bookmark around the location of the most recent\n// edit, so after changing the heading text we also relocate that bookmark\n// from wherever it used to sit (end of the document, from a prior edit) to\n// right after the newly-edited run.\n\n// 1) Drop the stale _GoBack bookmark (if present) before we touch anything.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the heading paragraph's run text and replace it.\nconst body = context.document.body;\nconst results = body.search(\"Weekly Progress Report #3 \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  // Insert the new text immediately before the old text so the new run\n  // inherits the same run formatting (font, size, etc.) as the original,\n  // then delete the old text, leaving only \"Weekly Progress Report #4\".\n  const newRange = target.insertText(\"Weekly Progress Report #4\", Word.InsertLocation.start);\n  await context.sync();\n\n  target.delete();\n  await context.sync();\n\n  // 3) Re-create the _GoBack bookmark immediately after the new text, i.e.\n  // collapsed to the end of the inserted range, matching Word's behavior of\n  // marking the point right after the last edit.\n  const endOfEdit = newRange.getRange(Word.RangeLocation.end);\n  endOfEdit.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word Weekly Progress Report heading: bump \"#3\" to \"#4\".\n# Word keeps a \"_GoBack\" bookmark around the location of the most recent\n# edit, so after changing the heading text we also relocate that bookmark\n# from wherever it used to sit (end of the document, from a prior edit) to\n# right after the newly-edited run.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the stale _GoBack bookmark (if present), wherever it currently sits.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Find the heading run that still says \"#3\" and replace its text with \"#4\".\n$target = $d.Content\n$found = $target.Find.Execute(\"Weekly Progress Report #3 \")\n\nif ($found) {\n    # Assigning .Text replaces the found range in place and collapses $target\n    # down to exactly the newly written text (\"Weekly Progress Report #4\").\n    $target.Text = \"Weekly Progress Report #4\"\n\n    # 3) Re-create the _GoBack bookmark immediately after the new text.\n    # A bookmark collapsed exactly on a run boundary can land in the wrong\n    # spot, so append a throwaway character right after the new text, anchor\n    # the bookmark just before it, then delete the throwaway character.\n    $tail = $target.Duplicate\n    $tail.Collapse(0)\n    $tail.InsertAfter(\"X\")\n\n    $bmPos = $tail.Start\n    $bmRange = $d.Range($bmPos, $bmPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n    $marker = $d.Range($bmPos, $bmPos + 1)\n    $marker.Delete()\n}\n"}
